# Apply the "#5: cash & deposit done" edit to the 存款 (deposits) sheet.
#
# The deposit sheet is being normalized to the same schema used by the
# other property sheets (stock, insurance, debt, ...):
#   bank | deposit_type | currency | owner | total | property_category |
#   category | date | legislator_name | legislator_id | source_file | index
#
# Concretely: the old "quantity" column (F, almost always blank) is
# dropped, the old "amount" column (G) becomes the new "total" column
# (F), and six new descriptive/provenance columns are appended
# (G..M): name/category, category, date, legislator_name,
# legislator_id, source_file, index.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)   # 存款

$headerStyleSrc = $ws.Range("B1")
$dataStyleSrc   = $ws.Range("B2")

# ---- Header row (row 1) --------------------------------------------
$ws.Range("E1").Value = "owner"
$ws.Range("F1").Value = "total"
$ws.Range("G1").Value = "property_category"
$ws.Range("H1").Value = "category"
$ws.Range("I1").Value = "date"
$ws.Range("J1").Value = "legislator_name"
$ws.Range("K1").Value = "legislator_id"
$ws.Range("L1").Value = "source_file"
$ws.Range("M1").Value = "index"

$headerStyleSrc.Copy()
$ws.Range("E1:M1").PasteSpecial(-4122)

# ---- Data rows (rows 2-14) ------------------------------------------
$lastRow = 14

# The "date" column holds a plain yyyy-mm-dd label (shared string in the
# source data), not a real Excel date value - force text format first so
# Excel does not auto-convert it into a date serial number.
$ws.Range("I2:I" + $lastRow).NumberFormat = "@"

for ($r = 2; $r -le $lastRow; $r++) {
    $idxVal = $ws.Cells.Item($r, 1).Value2          # column A: running index
    $amount = $ws.Cells.Item($r, 7).Value2           # old column G: amount

    $ws.Cells.Item($r, 6).Value2  = $amount           # F: total (was G)
    $ws.Cells.Item($r, 7).Value   = "deposit"         # G: property_category
    $ws.Cells.Item($r, 8).Value   = "normal"          # H: category
    $ws.Cells.Item($r, 9).Value   = "2011-11-22"      # I: date
    $ws.Cells.Item($r, 10).Value  = "丁守中"           # J: legislator_name
    $ws.Cells.Item($r, 11).Value2 = 515                # K: legislator_id
    $ws.Cells.Item($r, 12).Value  = "tmp8fef1"        # L: source_file
    $ws.Cells.Item($r, 13).Value2 = $idxVal            # M: index

    $dataStyleSrc.Copy()
    $ws.Range($ws.Cells.Item($r, 6), $ws.Cells.Item($r, 13)).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
